$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1668.7693
$ws.Range("I40").Value = 1042
$ws.Range("J40").Value = 2400
$ws.Range("K40").Value = 1042
$ws.Range("L40").Value = 2400
$ws.Range("M40").Value = -867
$ws.Range("N40").Value = -2750

$ws.Range("H125").Value = 3473.9375
$ws.Range("I125").Value = 3633.5
$ws.Range("J125").Value = 3420.75
$ws.Range("K125").Value = 32701.5
$ws.Range("L125").Value = 30786.75
$ws.Range("M125").Value = -30241.5
$ws.Range("N125").Value = -35706.75

$ws.Range("H137").Value = 1531.3513
$ws.Range("I137").Value = 1459.6842
$ws.Range("J137").Value = 1607
$ws.Range("K137").Value = 4379.0526
$ws.Range("L137").Value = 4821
$ws.Range("M137").Value = -1829.0526
$ws.Range("N137").Value = -9921

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1727.2632
$ws.Range("I2").Value = 1817.7273
$ws.Range("J2").Value = 1602.875
$ws.Range("K2").Value = 1817.7273
$ws.Range("L2").Value = 1602.875
$ws.Range("M2").Value = -1704.7273
$ws.Range("N2").Value = -1828.875

$ws.Range("H61").Value = 11560.333
$ws.Range("I61").Value = 1569.6154
$ws.Range("J61").Value = 76500
$ws.Range("K61").Value = 1569.6154
$ws.Range("L61").Value = 76500
$ws.Range("M61").Value = -1357.6154
$ws.Range("N61").Value = -76924

$ws.Range("H110").Value = 1151.25
$ws.Range("I110").Value = 1185.1538
$ws.Range("J110").Value = 1004.3333
$ws.Range("K110").Value = 1185.1538
$ws.Range("L110").Value = 1004.3333
$ws.Range("M110").Value = 859.8462
$ws.Range("N110").Value = -5094.3333

$ws.Range("H116").Value = 1727.2632
$ws.Range("I116").Value = 1817.7273
$ws.Range("J116").Value = 1602.875
$ws.Range("K116").Value = 1817.7273
$ws.Range("L116").Value = 1602.875
$ws.Range("M116").Value = 476.2727
$ws.Range("N116").Value = -6190.875

$ws.Range("H130").Value = 27571.428
$ws.Range("J130").Value = 27571.428
$ws.Range("L130").Value = 27571.428
$ws.Range("N130").Value = -37611.428

$ws.Range("H136").Value = 11560.333
$ws.Range("I136").Value = 1569.6154
$ws.Range("J136").Value = 76500
$ws.Range("K136").Value = 4708.8462
$ws.Range("L136").Value = 229500
$ws.Range("M136").Value = -2158.8462
$ws.Range("N136").Value = -234600

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1727.2632
$ws.Range("I3").Value = 1817.7273
$ws.Range("J3").Value = 1602.875
$ws.Range("K3").Value = 1817.7273
$ws.Range("L3").Value = 1602.875
$ws.Range("M3").Value = -1703.7273
$ws.Range("N3").Value = -1830.875

$ws.Range("H134").Value = 4226.6665
$ws.Range("I134").Value = 3898.0571
$ws.Range("K134").Value = 11694.1713
$ws.Range("M134").Value = -9159.1713

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2442.9333
$ws.Range("I134").Value = 2316.1904
$ws.Range("J134").Value = 2738.6667
$ws.Range("K134").Value = 6948.5712
$ws.Range("L134").Value = 8216.000100000001
$ws.Range("M134").Value = -4413.5712
$ws.Range("N134").Value = -13286.0001

$ws.Range("H141").Value = 56102.332
$ws.Range("J141").Value = 56102.332
$ws.Range("L141").Value = 56102.332
$ws.Range("N141").Value = -66462.33199999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 965.41
$ws.Range("I131").Value = 376.66666
$ws.Range("J131").Value = 983.6185
$ws.Range("K131").Value = 1129.99998
$ws.Range("L131").Value = 2950.8555
$ws.Range("M131").Value = 3910.00002
$ws.Range("N131").Value = -13030.8555

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2312.1177
$ws.Range("I80").Value = 2416.6667
$ws.Range("J80").Value = 2255.0908
$ws.Range("K80").Value = 2416.6667
$ws.Range("L80").Value = 2255.0908
$ws.Range("M80").Value = -1418.6667
$ws.Range("N80").Value = -4251.0908

$ws.Range("H83").Value = 2312.1177
$ws.Range("I83").Value = 2416.6667
$ws.Range("J83").Value = 2255.0908
$ws.Range("K83").Value = 12083.3335
$ws.Range("L83").Value = 11275.454
$ws.Range("M83").Value = -7091.333500000001
$ws.Range("N83").Value = -21259.454

$ws.Range("H113").Value = 1700
$ws.Range("I113").Value = 1700
$ws.Range("K113").Value = 1700
$ws.Range("M113").Value = 470

$ws.Range("H122").Value = 783.6875
$ws.Range("I122").Value = 763.6429000000001
$ws.Range("J122").Value = 924
$ws.Range("K122").Value = 2290.9287
$ws.Range("L122").Value = 2772
$ws.Range("M122").Value = 159.0712999999996
$ws.Range("N122").Value = -7672

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 59605.06
$ws.Range("I22").Value = 200546
$ws.Range("J22").Value = 879.6667
$ws.Range("K22").Value = 200546
$ws.Range("L22").Value = 879.6667
$ws.Range("M22").Value = -200251
$ws.Range("N22").Value = -1469.6667

$ws.Range("H27").Value = 59605.06
$ws.Range("I27").Value = 200546
$ws.Range("J27").Value = 879.6667
$ws.Range("K27").Value = 200546
$ws.Range("L27").Value = 879.6667
$ws.Range("M27").Value = -200439
$ws.Range("N27").Value = -1093.6667

$ws.Range("H40").Value = 2760.8948
$ws.Range("I40").Value = 2810.5334
$ws.Range("J40").Value = 2574.75
$ws.Range("K40").Value = 2810.5334
$ws.Range("L40").Value = 2574.75
$ws.Range("M40").Value = -2674.5334
$ws.Range("N40").Value = -2846.75

$ws.Range("H82").Value = 866.6667
$ws.Range("I82").Value = 620
$ws.Range("K82").Value = 620
$ws.Range("M82").Value = -259

$ws.Range("H85").Value = 866.6667
$ws.Range("I85").Value = 620
$ws.Range("K85").Value = 620
$ws.Range("M85").Value = 628

$ws.Range("H93").Value = 37862.684
$ws.Range("I93").Value = 1067.3572
$ws.Range("J93").Value = 140889.6
$ws.Range("K93").Value = 1067.3572
$ws.Range("L93").Value = 140889.6
$ws.Range("M93").Value = 180.6428000000001
$ws.Range("N93").Value = -143385.6

$ws.Range("H100").Value = 2000
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 2000
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 2000
$ws.Range("M100").ClearContents()
$ws.Range("N100").Value = -3082

$ws.Range("H136").Value = 1743.5676
$ws.Range("I136").Value = 1272.48
$ws.Range("J136").Value = 2725
$ws.Range("K136").Value = 3817.44
$ws.Range("L136").Value = 8175
$ws.Range("M136").Value = -1267.44
$ws.Range("N136").Value = -13275

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 600
$ws.Range("I96").Value = 600
$ws.Range("K96").Value = 600
$ws.Range("M96").Value = 773
